# Generate Report for Handoff
# Replace the old handoff id / hashes / timestamps with the new ones across
# all three worksheets (Overview, zh-cn, de-de), keeping the existing
# hyperlinks intact (only their displayed text changes).

$wb = $excel.ActiveWorkbook

$oldId = "c82872d3-7a9b-4c09-ae04-b95897554731"
$newId = "431f6735-b326-4f85-8d28-f0df21634fd2"

$oldHash = "421a6ffd6f34afe9c6ec138ef58a2564c95c38d9"
$newHash = "2cba08426a5dd0bdea3b32e17f0f2961cc4db448"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("D2").Value = "2016-45-12 12:45:07"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Range("D2").Value = "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-12 12:45:04"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Range("D2").Value = "$newId.$newHash.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-12 12:45:07"
